# Change "autonomous" to "connect" in the thesis description
# (commit: "Change autonomous to connected in thesis description")
$d = $word.ActiveDocument

$find = $d.Content.Find
$find.Text = "autonomous vehicles"
$find.Replacement.Text = "connect vehicles"
$find.Execute($find.Text, $true, $false, $false, $false, $false, `
              $true, 1, $false, $find.Replacement.Text, 2)
